$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.Value = "'52.511.86"
$c.Style = $s
$ws.Range("E2").Value = "  -13.22%  "
$c = $ws.Range("D3")
$s = $c.Style
$c.Value = "'2.305.79"
$c.Style = $s
$ws.Range("E3").Value = "  -20.55%  "
$ws.Range("E4").Value = "  +0.21%  "
$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'445.58"
$c.Style = $s
$ws.Range("E5").Value = "  -15.39%  "
$c = $ws.Range("D6")
$s = $c.Style
$c.Value = "'119.92"
$c.Style = $s
$ws.Range("E6").Value = "  -16.08%  "
$ws.Range("E7").Value = "  -0.17%  "
$c = $ws.Range("D8")
$s = $c.Style
$c.Value = "'0.466"
$c.Style = $s
$ws.Range("E8").Value = "  -15.19%  "
$c = $ws.Range("D9")
$s = $c.Style
$c.Value = "'2.316.73"
$c.Style = $s
$ws.Range("E9").Value = "  -20.27%  "
$c = $ws.Range("D10")
$s = $c.Style
$c.Value = "'5.28"
$c.Style = $s
$ws.Range("E10").Value = "  -11.46%  "
$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'0.0869"
$c.Style = $s
$ws.Range("E11").Value = "  -18.74%  "
$c = $ws.Range("D12")
$s = $c.Style
$c.Value = "'0.300"
$c.Style = $s
$ws.Range("E12").Value = "  -16.53%  "
$ws.Range("E13").Value = "  -5.95%  "
$c = $ws.Range("D14")
$s = $c.Style
$c.Value = "'52.507.70"
$c.Style = $s
$ws.Range("E14").Value = "  -13.23%  "
$c = $ws.Range("D15")
$s = $c.Style
$c.Value = "'18.72"
$c.Style = $s
$ws.Range("E15").Value = "  -17.28%  "
$ws.Range("E16").Value = "  -16.40%  "
$c = $ws.Range("D17")
$s = $c.Style
$c.Value = "'2.327.96"
$c.Style = $s
$ws.Range("E17").Value = "  -19.69%  "
$c = $ws.Range("D18")
$s = $c.Style
$c.Value = "'3.92"
$c.Style = $s
$ws.Range("E18").Value = "  -21.13%  "
$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'297.24"
$c.Style = $s
$ws.Range("E19").Value = "  -15.65%  "
$c = $ws.Range("D20")
$s = $c.Style
$c.Value = "'8.87"
$c.Style = $s
$ws.Range("E20").Value = "  -23.71%  "
$c = $ws.Range("D21")
$s = $c.Style
$c.Value = "'0.999"
$c.Style = $s
$ws.Range("E21").Value = "  -0.08%  "
$c = $ws.Range("D22")
$s = $c.Style
$c.Value = "'5.61"
$c.Style = $s
$ws.Range("E22").Value = "  -1.80%  "
$c = $ws.Range("D23")
$s = $c.Style
$c.Value = "'5.12"
$c.Style = $s
$ws.Range("E23").Value = "  -21.69%  "
$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'53.41"
$c.Style = $s
$ws.Range("E24").Value = "  -17.26%  "
$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'0.364"
$c.Style = $s
$ws.Range("E25").Value = "  -19.49%  "
$c = $ws.Range("D26")
$s = $c.Style
$c.Value = "'0.144"
$c.Style = $s
$ws.Range("E26").Value = "  -18.74%  "
$c = $ws.Range("D27")
$s = $c.Style
$c.Value = "'6.91"
$c.Style = $s
$ws.Range("E27").Value = "  -11.73%  "
$ws.Range("E28").Value = "  -0.27%  "
$c = $ws.Range("D29")
$s = $c.Style
$c.Value = "'0.0₃0658"
$c.Style = $s
$ws.Range("E29").Value = "  -22.28%  "
$c = $ws.Range("D30")
$s = $c.Style
$c.Value = "'140.75"
$c.Style = $s
$ws.Range("E30").Value = "  -6.36%  "
$c = $ws.Range("D31")
$s = $c.Style
$c.Value = "'16.77"
$c.Style = $s
$ws.Range("E31").Value = "  -14.37%  "
$c = $ws.Range("D32")
$s = $c.Style
$c.Value = "'1.34"
$c.Style = $s
$ws.Range("E32").Value = "  -19.93%  "
$c = $ws.Range("D33")
$s = $c.Style
$c.Value = "'4.71"
$c.Style = $s
$ws.Range("E33").Value = "  -15.23%  "
$c = $ws.Range("D34")
$s = $c.Style
$c.Value = "'0.823"
$c.Style = $s
$ws.Range("E34").Value = "  -17.50%  "
$c = $ws.Range("D35")
$s = $c.Style
$c.Value = "'3.39"
$c.Style = $s
$ws.Range("E35").Value = "  -21.65%  "
$c = $ws.Range("D36")
$s = $c.Style
$c.Value = "'0.996"
$c.Style = $s
$ws.Range("E36").Value = "  -0.05%  "
$c = $ws.Range("D37")
$s = $c.Style
$c.Value = "'0.990"
$c.Style = $s
$ws.Range("E37").Value = "  -17.38%  "
$c = $ws.Range("D38")
$s = $c.Style
$c.Value = "'31.80"
$c.Style = $s
$ws.Range("E38").Value = "  -15.73%  "
$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'10.15"
$c.Style = $s
$ws.Range("E39").Value = "  -1.63%  "
$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'0.556"
$c.Style = $s
$ws.Range("E40").Value = "  -14.20%  "
$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'0.0503"
$c.Style = $s
$ws.Range("E41").Value = "  -13.33%  "
$c = $ws.Range("D42")
$s = $c.Style
$c.Value = "'3.09"
$c.Style = $s
$ws.Range("E42").Value = "  -16.54%  "
$c = $ws.Range("D43")
$s = $c.Style
$c.Value = "'1.909.16"
$c.Style = $s
$ws.Range("E43").Value = "  -16.54%  "
$c = $ws.Range("D44")
$s = $c.Style
$c.Value = "'1.15"
$c.Style = $s
$ws.Range("E44").Value = "  -21.23%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'0.0821"
$c.Style = $s
$ws.Range("E45").Value = "  -10.98%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'4.18"
$c.Style = $s
$ws.Range("E46").Value = "  -15.56%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D47")
$s = $c.Style
$c.Value = "'0.0204"
$c.Style = $s
$ws.Range("E47").Value = "  -14.02%  "
$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'15.47"
$c.Style = $s
$ws.Range("E48").Value = "  -24.37%  "
$ws.Range("E49").Value = "  -5.24%  "
$c = $ws.Range("D50")
$s = $c.Style
$c.Value = "'4.42"
$c.Style = $s
$ws.Range("E50").Value = "  -13.46%  "
$ws.Range("E51").Value = "  -18.39%  "
